$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "56.697.12"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
Set-TextValue "D3" "2.408.61"
$ws.Range("E3").Value = "  -3.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
Set-TextValue "D5" "486.68"
$ws.Range("E5").Value = "  -1.76%  "

# Row 6
Set-TextValue "D6" "152.89"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
Set-TextValue "D8" "0.603"
$ws.Range("E8").Value = "  +17.34%  "

# Row 9
Set-TextValue "D9" "2.426.14"
$ws.Range("E9").Value = "  -3.61%  "

# Row 10
Set-TextValue "D10" "0.0996"
$ws.Range("E10").Value = "  +0.29%  "

# Row 11
Set-TextValue "D11" "5.73"
$ws.Range("E11").Value = "  -1.02%  "

# Row 12
Set-TextValue "D12" "0.335"
$ws.Range("E12").Value = "  -0.58%  "

# Row 13
$ws.Range("E13").Value = "  +1.37%  "

# Row 14
Set-TextValue "D14" "2.828.15"
$ws.Range("E14").Value = "  -3.59%  "

# Row 15
Set-TextValue "D15" "56.918.89"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16
Set-TextValue "D16" "20.74"
$ws.Range("E16").Value = "  -3.66%  "

# Row 17
$ws.Range("E17").Value = "  -2.82%  "

# Row 18
Set-TextValue "D18" "2.421.54"
$ws.Range("E18").Value = "  -4.00%  "

# Row 19
$ws.Range("E19").Value = "  +3.78%  "

# Row 20
Set-TextValue "D20" "324.26"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("E21").Value = "  -4.18%  "

# Row 22
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23
Set-TextValue "D23" "5.94"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
Set-TextValue "D24" "57.93"
$ws.Range("E24").Value = "  -1.97%  "

# Row 25
Set-TextValue "D25" "0.408"
$ws.Range("E25").Value = "  -1.02%  "

# Row 26
Set-TextValue "D26" "0.996"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
Set-TextValue "D27" "0.160"
$ws.Range("E27").Value = "  -2.69%  "

# Row 28
Set-TextValue "D28" "2.516.66"
$ws.Range("E28").Value = "  -3.76%  "

# Row 29
Set-TextValue "D29" "7.30"
$ws.Range("E29").Value = "  -4.77%  "

# Row 30
Set-TextValue "D30" "0.0₃0785"
$ws.Range("E30").Value = "  -4.73%  "

# Row 31
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "149.46"
$ws.Range("E32").Value = "  -1.95%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "18.58"
$ws.Range("E33").Value = "  +0.78%  "

# Row 34
Set-TextValue "D34" "1.52"
$ws.Range("E34").Value = "  -0.89%  "

# Row 35
$ws.Range("E35").Value = "  +1.59%  "

# Row 36
Set-TextValue "D36" "1.15"
$ws.Range("E36").Value = "  -2.08%  "

# Row 37
Set-TextValue "D37" "3.71"
$ws.Range("E37").Value = "  -2.33%  "

# Row 38
Set-TextValue "D38" "0.845"
$ws.Range("E38").Value = "  -3.54%  "

# Row 39
$ws.Range("E39").Value = "  +9.13%  "

# Row 40
Set-TextValue "D40" "34.11"
$ws.Range("E40").Value = "  -0.51%  "

# Row 41
Set-TextValue "D41" "3.53"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("E42").Value = "  -2.17%  "

# Row 43
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D44" "269.77"
$ws.Range("E44").Value = "  +0.76%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.591"
$ws.Range("E45").Value = "  -4.05%  "

# Row 46
$ws.Range("E46").Value = "  -6.15%  "

# Row 47
Set-TextValue "D47" "10.20"
$ws.Range("E47").Value = "  -0.16%  "

# Row 48
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
Set-TextValue "D49" "4.62"
$ws.Range("E49").Value = "  -6.72%  "

# Row 50
$ws.Range("E50").Value = "  -3.17%  "

# Row 51
Set-TextValue "D51" "1.865.20"
$ws.Range("E51").Value = "  -2.48%  "
